$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 16 (shifts the old "total" row from 16 down to 17)
$ws.Rows.Item(16).Insert()

$ws.Range("A16").Value = "P148"
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1.7

$ws.Range("B17").Select()
